$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated values for columns B, C, E, F across rows 2-15 (D column unchanged)
$data = @{
    2  = @{ B = 7623.24506588005;  C = 7264.34498419236;  E = 3084.94424563731;  F = -20.6129487570968 }
    3  = @{ B = 7464.92848859504;  C = 6943.15174689051;  E = 2811.24334987376;  F = 255.599795698511 }
    4  = @{ B = 6989.25126626053;  C = 5873.20277713754;  E = 3271.3184751121;   F = 230.188385510402 }
    5  = @{ B = 2155.54830302928;  C = 3456.37655075567;  E = 3297.35409307157;  F = -193.761223173865 }
    6  = @{ B = 2116.37887855762;  C = 3655.52633817196;  E = 3257.67647985881;  F = -187.116549248718 }
    7  = @{ B = 7968.62117480238;  C = 6949.41607719639;  E = 4215.55112582053;  F = -9.95969987429506 }
    8  = @{ B = 8517.20594041771;  C = 7639.30403838576;  E = 4596.07156508061;  F = 34.6406501444323 }
    9  = @{ B = 8517.20594041771;  C = 7791.35957696154;  E = 4596.07156508061;  F = 40.9762975850896 }
    10 = @{ B = 8525.29585685903;  C = 7968.83765599431;  E = 4598.02789310823;  F = 48.4527312126058 }
    11 = @{ B = 8530.22044936565;  C = 7178.21516842877;  E = 4628.66379180446;  F = 16.7866233430514 }
    12 = @{ B = 3118.33163974066;  C = 4880.64410128015;  E = 4317.33067620137;  F = -91.917717604937 }
    13 = @{ B = 2984.2957163966;   C = 4801.50987573279;  E = 4310.06456827996;  F = -95.5177314994689 }
    14 = @{ B = 9278.30282674808;  C = 8113.97122708144;  E = 5257.97981428584;  F = 81.9979600569701 }
    15 = @{ B = 9278.30412871915;  C = 8610.89742141984;  E = 5257.98357839086;  F = 102.703374992113 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("F$row").Value = $vals.F
}
